$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2318840579710145
$ws.Range("C2").Value = 0.4855072463768116
$ws.Range("J2").Value = 0.01159420289855072
$ws.Range("P2").Value = 0.1797101449275362
$ws.Range("S2").Value = 0.09130434782608696
$ws.Range("B3").Value = 0.01432664756446991
$ws.Range("C3").Value = 0.0171919770773639
$ws.Range("J3").Value = 0.02005730659025788
$ws.Range("P3").Value = 0.7621776504297995
$ws.Range("S3").Value = 0.1862464183381089
$ws.Range("J4").Value = 0.02040816326530612
$ws.Range("P4").Value = 0.7142857142857143
$ws.Range("S4").Value = 0.2653061224489796
$ws.Range("B6").Value = 0.08418891170431211
$ws.Range("D6").Value = 0.008213552361396304
$ws.Range("E6").Value = 0.002053388090349076
$ws.Range("F6").Value = 0.07392197125256673
$ws.Range("J6").Value = 0.2012320328542095
$ws.Range("O6").Value = 0.02053388090349076
$ws.Range("Q6").Value = 0.1909650924024641
$ws.Range("R6").Value = 0.05749486652977413
$ws.Range("S6").Value = 0.3613963039014374
$ws.Range("B7").Value = 0.1216216216216216
$ws.Range("D7").Value = 0.02162162162162162
$ws.Range("E7").Value = 0.002702702702702703
$ws.Range("F7").Value = 0.02972972972972973
$ws.Range("J7").Value = 0.1648648648648649
$ws.Range("O7").Value = 0.03513513513513514
$ws.Range("Q7").Value = 0.2027027027027027
$ws.Range("R7").Value = 0.04594594594594595
$ws.Range("S7").Value = 0.3756756756756757
$ws.Range("B8").Value = 0.0906183368869936
$ws.Range("D8").Value = 0.0255863539445629
$ws.Range("F8").Value = 0.06716417910447761
$ws.Range("J8").Value = 0.1226012793176972
$ws.Range("O8").Value = 0.02452025586353945
$ws.Range("Q8").Value = 0.2046908315565032
$ws.Range("R8").Value = 0.07142857142857142
$ws.Range("S8").Value = 0.3933901918976546
$ws.Range("B9").Value = 0.08775981524249422
$ws.Range("D9").Value = 0.009237875288683603
$ws.Range("F9").Value = 0.07852193995381063
$ws.Range("J9").Value = 0.1177829099307159
$ws.Range("O9").Value = 0.02540415704387991
$ws.Range("Q9").Value = 0.2078521939953811
$ws.Range("R9").Value = 0.07621247113163972
$ws.Range("S9").Value = 0.3972286374133949
$ws.Range("B10").Value = 0.1125226860254084
$ws.Range("D10").Value = 0.02286751361161524
$ws.Range("F10").Value = 0.07005444646098004
$ws.Range("J10").Value = 0.1168784029038113
$ws.Range("O10").Value = 0.02032667876588022
$ws.Range("Q10").Value = 0.2246823956442831
$ws.Range("R10").Value = 0.06678765880217785
$ws.Range("S10").Value = 0.3658802177858439
$ws.Range("G11").Value = 0.1516666666666667
$ws.Range("J11").Value = 0.1033333333333333
$ws.Range("K11").Value = 0.2083333333333333
$ws.Range("L11").Value = 0.5216666666666666
$ws.Range("S11").Value = 0.015
$ws.Range("G12").Value = 0.7586206896551724
$ws.Range("J12").Value = 0.1724137931034483
$ws.Range("K12").Value = 0.01880877742946709
$ws.Range("L12").Value = 0.01880877742946709
$ws.Range("S12").Value = 0.03134796238244514
$ws.Range("G13").Value = 0.6025641025641025
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("S13").Value = 0.0641025641025641
$ws.Range("F15").Value = 0.02448979591836735
$ws.Range("H15").Value = 0.1489795918367347
$ws.Range("I15").Value = 0.06734693877551021
$ws.Range("J15").Value = 0.3836734693877551
$ws.Range("K15").Value = 0.05510204081632653
$ws.Range("M15").Value = 0.01020408163265306
$ws.Range("O15").Value = 0.05510204081632653
$ws.Range("S15").Value = 0.2551020408163265
$ws.Range("F16").Value = 0.03146067415730337
$ws.Range("H16").Value = 0.1797752808988764
$ws.Range("I16").Value = 0.07640449438202247
$ws.Range("J16").Value = 0.4292134831460674
$ws.Range("K16").Value = 0.1123595505617977
$ws.Range("M16").Value = 0.006741573033707865
$ws.Range("N16").Value = 0.002247191011235955
$ws.Range("O16").Value = 0.05168539325842696
$ws.Range("S16").Value = 0.1101123595505618
$ws.Range("F17").Value = 0.01785714285714286
$ws.Range("H17").Value = 0.1691729323308271
$ws.Range("I17").Value = 0.08834586466165413
$ws.Range("J17").Value = 0.443609022556391
$ws.Range("K17").Value = 0.07706766917293233
$ws.Range("M17").Value = 0.01503759398496241
$ws.Range("N17").Value = 0.001879699248120301
$ws.Range("O17").Value = 0.06203007518796992
$ws.Range("S17").Value = 0.125
$ws.Range("F18").Value = 0.0303951367781155
$ws.Range("H18").Value = 0.2036474164133739
$ws.Range("I18").Value = 0.1185410334346505
$ws.Range("J18").Value = 0.3799392097264438
$ws.Range("K18").Value = 0.08206686930091185
$ws.Range("M18").Value = 0.00911854103343465
$ws.Range("O18").Value = 0.06382978723404255
$ws.Range("S18").Value = 0.1124620060790274
$ws.Range("F19").Value = 0.01511715797430083
$ws.Range("H19").Value = 0.2055933484504913
$ws.Range("I19").Value = 0.08843537414965986
$ws.Range("J19").Value = 0.382464096749811
$ws.Range("K19").Value = 0.1016628873771731
$ws.Range("M19").Value = 0.02040816326530612
$ws.Range("N19").Value = 0.001511715797430083
$ws.Range("O19").Value = 0.07180650037792895
$ws.Range("S19").Value = 0.1130007558578987
